$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 219-220, pushing the existing data (old rows 219-317)
# down to become rows 221-319.
$ws.Rows("219:220").Insert()

# Populate the first new row (row 219)
$ws.Range("A219").Value = 11
$ws.Range("B219").Value = "Vega Monumental Concepción"
$ws.Range("C219").Value = "Bíobío"
$ws.Range("D219").Value = 44924
$ws.Range("E219").Value = 8
$ws.Range("F219").Value = 100114013
$ws.Range("G219").Value = "Zanahoria"
$ws.Range("H219").Value = "Sin especificar"
$ws.Range("I219").Value = "Primera"
$ws.Range("J219").Value = 800
$ws.Range("K219").Value = 7000
$ws.Range("L219").Value = 8000
$ws.Range("M219").Value = 7500
$ws.Range("N219").Value = "$/saco 20 kilos"
$ws.Range("O219").Value = "Región de Ñuble"
$ws.Range("P219").Value = 375
$ws.Range("Q219").Value = 20
$ws.Range("R219").Value = "Hortaliza"

# Populate the second new row (row 220)
$ws.Range("A220").Value = 11
$ws.Range("B220").Value = "Vega Monumental Concepción"
$ws.Range("C220").Value = "Bíobío"
$ws.Range("D220").Value = 44924
$ws.Range("E220").Value = 8
$ws.Range("F220").Value = 100114013
$ws.Range("G220").Value = "Zanahoria"
$ws.Range("H220").Value = "Sin especificar"
$ws.Range("I220").Value = "Segunda"
$ws.Range("J220").Value = 400
$ws.Range("K220").Value = 6000
$ws.Range("L220").Value = 6000
$ws.Range("M220").Value = 6000
$ws.Range("N220").Value = "$/saco 20 kilos"
$ws.Range("O220").Value = "Región de Ñuble"
$ws.Range("P220").Value = 300
$ws.Range("Q220").Value = 20
$ws.Range("R220").Value = "Hortaliza"
